$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.578.25"
$ws.Range("E2").Value = "  -5.25%  "
$ws.Range("D3").Value = "'3.060.46"
$ws.Range("E3").Value = "  -5.44%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'539.19"
$ws.Range("E5").Value = "  -7.25%  "
$ws.Range("D6").Value = "'132.91"
$ws.Range("E6").Value = "  -12.26%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'3.050.79"
$ws.Range("E8").Value = "  -5.34%  "
$ws.Range("E9").Value = "  -4.92%  "
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = "  -5.79%  "
$ws.Range("D11").Value = "'6.15"
$ws.Range("E11").Value = "  -13.32%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  -5.43%  "
$ws.Range("D13").Value = "'34.61"
$ws.Range("E13").Value = "  -8.03%  "
$ws.Range("E14").Value = "  -6.17%  "
$ws.Range("D15").Value = "'3.515.26"
$ws.Range("E15").Value = "  -6.44%  "
$ws.Range("D16").Value = "'62.585.12"
$ws.Range("E16").Value = "  -5.39%  "
$ws.Range("E17").Value = "  -2.91%  "
$ws.Range("D18").Value = "'3.061.20"
$ws.Range("E18").Value = "  -5.70%  "
$ws.Range("D19").Value = "'6.61"
$ws.Range("E19").Value = "  -6.74%  "
$ws.Range("D20").Value = "'478.84"
$ws.Range("E20").Value = "  -12.38%  "
$ws.Range("D21").Value = "'13.36"
$ws.Range("E21").Value = "  -7.77%  "
$ws.Range("D22").Value = "'0.704"
$ws.Range("E22").Value = "  -5.04%  "
$ws.Range("D23").Value = "'7.20"
$ws.Range("E23").Value = "  -7.99%  "
$ws.Range("D24").Value = "'78.42"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").Value = "'12.06"
$ws.Range("E25").Value = "  -10.18%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -8.86%  "
$ws.Range("D28").Value = "'8.19"
$ws.Range("E28").Value = "  -12.01%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'1.91"
$ws.Range("E30").Value = "  -14.48%  "
$ws.Range("D31").Value = "'25.99"
$ws.Range("E31").Value = "  -5.98%  "
$ws.Range("D32").Value = "'1.09"
$ws.Range("E32").Value = "  -6.52%  "
$ws.Range("D33").Value = "'2.42"
$ws.Range("E33").Value = "  -11.94%  "
$ws.Range("D34").Value = "'58.23"
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.97"
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "'483.85"
$ws.Range("E36").Value = "  -14.77%  "
$ws.Range("D37").Value = "'5.17"
$ws.Range("E37").Value = "  -8.43%  "
$ws.Range("D38").Value = "'3.133.47"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("D39").Value = "'0.0390"
$ws.Range("E39").Value = "  -13.73%  "
$ws.Range("D40").Value = "'0.0793"
$ws.Range("E40").Value = "  -7.58%  "
$ws.Range("E41").Value = "  -10.26%  "
$ws.Range("E42").Value = "  -6.30%  "
$ws.Range("D43").Value = "'2.55"
$ws.Range("E43").Value = "  -13.51%  "
$ws.Range("E44").Value = "  -10.50%  "
$ws.Range("D46").Value = "'2.04"
$ws.Range("E46").Value = "  -11.11%  "
$ws.Range("D47").Value = "'24.47"
$ws.Range("E47").Value = "  -7.20%  "
$ws.Range("D48").Value = "'118.45"
$ws.Range("E48").Value = "  -5.54%  "
$ws.Range("E49").Value = "  -4.73%  "
$ws.Range("D50").Value = "'0.0₃0505"
$ws.Range("E50").Value = "  -9.22%  "
$ws.Range("E51").Value = "  -8.82%  "
